# multiple_session_batch.xlsx batch-definition update
#
# - Insert two new parameter rows right after row 17 (the "Fuel Scenario
#   Annual Data File" row), pushing every row below down by two.
# - Row 17 becomes "Context Fuel Prices File" (renamed in place).
# - New row 18 becomes "Context New Vehicle Market File" pointing at the
#   new context_new_vehicle_market.csv input file.
# - New row 19 keeps the original "Fuel Scenario Annual Data File" label,
#   still pointing at context_fuel_prices.csv.
# Everything else on the sheet shifts down by two rows automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for two new rows at 18 and 19, copying formatting/content from
# row 17 so the inserted rows inherit its styling, then fix up the values.
$ws.Range("A18:S19").Insert()
$ws.Range("A17:S17").Copy($ws.Range("A18:S18"))
$ws.Range("A17:S17").Copy($ws.Range("A19:S19"))

# Row 17: "Fuel Scenario Annual Data File" -> "Context Fuel Prices File"
# (value/C/D columns stay pointed at context_fuel_prices.csv)
$ws.Range("A17").Value = "Context Fuel Prices File"

# Row 18 (new): Context New Vehicle Market File
$ws.Range("A18").Value = "Context New Vehicle Market File"
$ws.Range("B18").Value = "String"
$ws.Range("C18").Value = "input_samples/context_new_vehicle_market.csv"
$ws.Range("D18").Value = "input_samples/context_new_vehicle_market.csv"

# Row 19 (new): keeps the old "Fuel Scenario Annual Data File" label,
# still referencing context_fuel_prices.csv (copied from row 17 already).
$ws.Range("A19").Value = "Fuel Scenario Annual Data File"
$ws.Range("B19").Value = "String"
$ws.Range("C19").Value = "input_samples/context_fuel_prices.csv"
$ws.Range("D19").Value = "input_samples/context_fuel_prices.csv"
